$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 374.66666
$ws.Range("I38").Value = 196
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 588
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -216
$ws.Range("N38").Value = -3744
$ws.Range("H43").Value = 874.6
$ws.Range("I43").Value = 800.5
$ws.Range("J43").Value = 893.125
$ws.Range("K43").Value = 800.5
$ws.Range("L43").Value = 893.125
$ws.Range("M43").Value = -731.5
$ws.Range("N43").Value = -1031.125
$ws.Range("H116").Value = 10701325
$ws.Range("I116").Value = 21392998
$ws.Range("J116").Value = 9650.833000000001
$ws.Range("K116").Value = 21392998
$ws.Range("L116").Value = 9650.833000000001
$ws.Range("M116").Value = -21389556
$ws.Range("N116").Value = -16534.833
$ws.Range("H132").Value = 275017.2
$ws.Range("I132").Value = 310575.28
$ws.Range("J132").Value = 61668.668
$ws.Range("K132").Value = 931725.8400000001
$ws.Range("L132").Value = 185006.004
$ws.Range("M132").Value = -929195.8400000001
$ws.Range("N132").Value = -190066.004
$ws.Range("H138").Value = 12123827
$ws.Range("I138").Value = 5556527
$ws.Range("K138").Value = 16669581
$ws.Range("M138").Value = -16664441
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 100001
$ws.Range("I57").Value = 100001
$ws.Range("K57").Value = 100001
$ws.Range("M57").Value = -99517
$ws.Range("H74").Value = 7454.6665
$ws.Range("I74").Value = 1863.2
$ws.Range("J74").Value = 21433.334
$ws.Range("K74").Value = 1863.2
$ws.Range("L74").Value = 21433.334
$ws.Range("M74").Value = -989.2
$ws.Range("N74").Value = -23181.334
$ws.Range("H77").Value = 7454.6665
$ws.Range("I77").Value = 1863.2
$ws.Range("J77").Value = 21433.334
$ws.Range("K77").Value = 9316
$ws.Range("L77").Value = 107166.67
$ws.Range("M77").Value = -4948
$ws.Range("N77").Value = -115902.67
$ws.Range("H110").Value = 1005.61536
$ws.Range("I110").Value = 855
$ws.Range("J110").Value = 1507.6666
$ws.Range("K110").Value = 855
$ws.Range("L110").Value = 1507.6666
$ws.Range("M110").Value = 1190
$ws.Range("N110").Value = -5597.6666
$ws.Range("H132").Value = 3538.0833
$ws.Range("I132").Value = 3340
$ws.Range("J132").Value = 4132.3335
$ws.Range("K132").Value = 10020
$ws.Range("L132").Value = 12397.0005
$ws.Range("M132").Value = -7490
$ws.Range("N132").Value = -17457.0005
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 56101.11
$ws.Range("I16").Value = 71773
$ws.Range("K16").Value = 71773
$ws.Range("M16").Value = -71486
$ws.Range("H31").Value = 1197.6072
$ws.Range("I31").Value = 1089.32
$ws.Range("J31").Value = 2100
$ws.Range("K31").Value = 1089.32
$ws.Range("L31").Value = 2100
$ws.Range("M31").Value = -794.3199999999999
$ws.Range("N31").Value = -2690
$ws.Range("H34").Value = 1197.6072
$ws.Range("I34").Value = 1089.32
$ws.Range("J34").Value = 2100
$ws.Range("K34").Value = 1089.32
$ws.Range("L34").Value = 2100
$ws.Range("M34").Value = -887.3199999999999
$ws.Range("N34").Value = -2504
$ws.Range("H58").Value = 2106.7083
$ws.Range("J58").Value = 3536.125
$ws.Range("L58").Value = 3536.125
$ws.Range("N58").Value = -3942.125
$ws.Range("H99").Value = 7813700
$ws.Range("I99").Value = 10417600
$ws.Range("J99").Value = 1999.5
$ws.Range("K99").Value = 10417600
$ws.Range("L99").Value = 1999.5
$ws.Range("M99").Value = -10416102
$ws.Range("N99").Value = -4995.5
$ws.Range("H105").Value = 735.9091
$ws.Range("I105").Value = 609.6
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 609.6
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = 1137.4
$ws.Range("N105").Value = -5493
$ws.Range("H113").Value = 56101.11
$ws.Range("I113").Value = 71773
$ws.Range("K113").Value = 71773
$ws.Range("M113").Value = -69603
$ws.Range("H122").Value = 2086.8
$ws.Range("I122").Value = 1432.1
$ws.Range("K122").Value = 4296.299999999999
$ws.Range("M122").Value = -1846.299999999999
$ws.Range("H126").Value = 7813700
$ws.Range("I126").Value = 10417600
$ws.Range("J126").Value = 1999.5
$ws.Range("K126").Value = 31252800
$ws.Range("L126").Value = 5998.5
$ws.Range("M126").Value = -31250330
$ws.Range("N126").Value = -10938.5
$ws.Range("H136").Value = 2106.7083
$ws.Range("J136").Value = 3536.125
$ws.Range("L136").Value = 10608.375
$ws.Range("N136").Value = -15708.375
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1284.6538
$ws.Range("I5").Value = 979
$ws.Range("J5").Value = 2303.5
$ws.Range("K5").Value = 2937
$ws.Range("L5").Value = 6910.5
$ws.Range("M5").Value = -2825
$ws.Range("N5").Value = -7134.5
$ws.Range("H131").Value = 1772.6757
$ws.Range("I131").Value = 570
$ws.Range("J131").Value = 2005.4517
$ws.Range("K131").Value = 1710
$ws.Range("L131").Value = 6016.355100000001
$ws.Range("M131").Value = 3330
$ws.Range("N131").Value = -16096.3551
$ws.Range("H132").Value = 1436.0625
$ws.Range("J132").Value = 1465.1333
$ws.Range("L132").Value = 13186.1997
$ws.Range("N132").Value = -18246.1997
$ws.Range("H133").Value = 6058.8887
$ws.Range("J133").Value = 11750
$ws.Range("L133").Value = 35250
$ws.Range("N133").Value = -45370
$ws.Range("H135").Value = 1284.6538
$ws.Range("I135").Value = 979
$ws.Range("J135").Value = 2303.5
$ws.Range("K135").Value = 8811
$ws.Range("L135").Value = 20731.5
$ws.Range("M135").Value = -6276
$ws.Range("N135").Value = -25801.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1000000
$ws.Range("I97").Value = 1000000
$ws.Range("K97").Value = 1000000
$ws.Range("M97").Value = -999504
$ws.Range("H113").Value = 2499
$ws.Range("I113").Value = 2499
$ws.Range("K113").Value = 2499
$ws.Range("M113").Value = -329
$ws.Range("H122").Value = 741610.0600000001
$ws.Range("I122").Value = 1111614
$ws.Range("J122").Value = 1602.2
$ws.Range("K122").Value = 3334842
$ws.Range("L122").Value = 4806.6
$ws.Range("M122").Value = -3332392
$ws.Range("N122").Value = -9706.6
$ws.Range("H132").Value = 3890
$ws.Range("I132").Value = 1972.5
$ws.Range("J132").Value = 5807.5
$ws.Range("K132").Value = 5917.5
$ws.Range("L132").Value = 17422.5
$ws.Range("M132").Value = -3387.5
$ws.Range("N132").Value = -22482.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3249.1738
$ws.Range("I122").Value = 1771.2858
$ws.Range("J122").Value = 3895.75
$ws.Range("K122").Value = 5313.857400000001
$ws.Range("L122").Value = 11687.25
$ws.Range("N122").Value = -16587.25
$ws.Range("M122").Value = -2863.857400000001
$ws.Range("H132").Value = 4142.5415
$ws.Range("I132").Value = 3050.6667
$ws.Range("J132").Value = 5234.4165
$ws.Range("K132").Value = 9152.000100000001
$ws.Range("L132").Value = 15703.2495
$ws.Range("M132").Value = -6622.000100000001
$ws.Range("N132").Value = -20763.2495
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3906.44
$ws.Range("I81").Value = 1110.1666
$ws.Range("J81").Value = 4789.4736
$ws.Range("K81").Value = 2220.3332
$ws.Range("L81").Value = 9578.947200000001
$ws.Range("M81").Value = -1159.3332
$ws.Range("N81").Value = -11700.9472
$ws.Range("H84").Value = 3906.44
$ws.Range("I84").Value = 1110.1666
$ws.Range("J84").Value = 4789.4736
$ws.Range("K84").Value = 11101.666
$ws.Range("L84").Value = 47894.736
$ws.Range("M84").Value = -5797.666000000001
$ws.Range("N84").Value = -58502.736
$ws.Range("H113").Value = 163.33333
$ws.Range("I113").Value = 163.33333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 489.99999
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1680.00001
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 49195.668
$ws.Range("I122").Value = 72672
$ws.Range("J122").Value = 2243
$ws.Range("K122").Value = 218016
$ws.Range("L122").Value = 6729
$ws.Range("M122").Value = -215566
$ws.Range("N122").Value = -11629
$ws.Range("H132").Value = 12502831
$ws.Range("I132").Value = 16668680
$ws.Range("J132").Value = 5284
$ws.Range("K132").Value = 50006040
$ws.Range("L132").Value = 15852
$ws.Range("M132").Value = -50003510
$ws.Range("N132").Value = -20912
